$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# New run entry added to row 24 (map6 / map11 / "all" flagged as "R", plus run dir note)
$ws.Range("H24").Value = "R"
$ws.Range("M24").Value = "R"
$ws.Range("N24").Value = "R"
$ws.Range("O24").Value = "map*_use_dist_self_play_shuffle_pos_[actor]2layers_[frames]3"

# Update window view: position/size, zoom level, scroll position, and current selection
$excel.ActiveWindow.Left = 30860
$excel.ActiveWindow.Top = 1640
$excel.ActiveWindow.Width = 36200
$excel.ActiveWindow.Height = 18120
$excel.ActiveWindow.Zoom = 133
$excel.ActiveWindow.ScrollColumn = 7
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("O28").Select()
